$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.253.35"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.856.06"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'314.05"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.4661"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").Value = "'0.3708"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "'0.8906"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").Value = "'20.08"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").Value = "'0.07873"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "1.828.72"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "'5.409"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'6.514"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "'91.72"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "'0.000008922"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "27.289.58"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "'5.079"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'10.53"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "2.086.76"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").Value = "'2.030"
$ws.Range("E25").Value = "  +9.63%  "
$ws.Range("D26").Value = "'151.78"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'18.46"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "'2.045"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'115.98"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'3.148"
$ws.Range("E32").Value = "  +6.50%  "
$ws.Range("D33").Value = "'0.7694"
$ws.Range("E33").Value = "  +5.36%  "
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("D35").Value = "'4.527"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").Value = "'2.720"
$ws.Range("E36").Value = "  +10.42%  "
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("D38").Value = "'0.01943"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "'0.05234"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "'2.946"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'7.071"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'0.5121"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'0.1628"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").Value = "'8.534"
$ws.Range("E44").Value = "  +4.97%  "
$ws.Range("D45").Value = "'0.4793"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'10.37"
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'102.88"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").Value = "'1.646"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").Value = "'0.06202"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'65.55"
$ws.Range("E51").Value = "  +1.69%  "
